# Append the new Argent price row (row 83) to the "Prices" sheet, mirroring
# the existing rows where every column is stored as plain text (the sheet
# has no numeric formatting anywhere). A leading apostrophe forces Excel to
# keep each value as literal text instead of inferring a date/number type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A83").Value = "'2025-05-23"
$ws.Range("B83").Value = "'35.5"
$ws.Range("C83").Value = "'35.4"
$ws.Range("D83").Value = "'0.94"
$ws.Range("E83").Value = "'0.258"
$ws.Range("F83").Value = "'0.09"
$ws.Range("G83").Value = "'5,406"
$ws.Range("H83").Value = "'8,094"
$ws.Range("I83").Value = "'8,144"
$ws.Range("J83").Value = "'7.2186"
